$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010 and 2010-18")

# Insert a new row at 89; everything below (old rows 89-102) shifts down to 90-103,
# and Excel auto-adjusts all formula references accordingly.
$ws.Rows.Item(89).Insert()

# Populate the new row 89 with data (copy formatting from row 88, then set values).
$ws.Range("A88:R88").Copy()
$ws.Range("A89:R89").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A89").Value = "CW3M"
$ws.Range("B89").Value = "Demo_Baseline_2010-18_C474+"
$ws.Range("C89").Value = "2010-18"
$ws.Range("D89").Value = 929.71728533333328
$ws.Range("E89").Value = 1890.2624918888889
$ws.Range("F89").Value = 1.0680259999999999
$ws.Range("G89").Value = 270.41205844444437
$ws.Range("H89").Value = 9.8445367777777779
$ws.Range("I89").Value = 7.7082955555555559
$ws.Range("J89").Value = 8.2027718888888881
$ws.Range("K89").Value = 669.04810911111099
$ws.Range("L89").Value = 80.5032391111111
$ws.Range("M89").Value = 1418.8638372222224
$ws.Range("N89").Value = 932.63113755555548
$ws.Range("O89").Value = 5820.4378255555557
$ws.Range("P89").Value = 27412.728515555555
$ws.Range("Q89").Value = 0.23640099999999997
$ws.Range("R89").Value = [double]"4.7777777777777777E-5"

# remove highlight fill on M89/Q89/R89 (row 88 has highlighted Q/R/M style; new row 89 should not)
$ws.Range("M89").Interior.ColorIndex = -4142
$ws.Range("Q89").Interior.ColorIndex = -4142
$ws.Range("R89").Interior.ColorIndex = -4142

$ws.Range("V90").Select()
